$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)

# Remove the empty "TextBox 24" shape (id=25) that only contained an
# empty paragraph with no visible text.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 24") {
        $sh.Delete()
    }
}
